# Updates the cryptocurrency price/volume table (and the swapped
# Quant / EnergySwap row pair) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.332.72'
$ws.Range("E2").Value = '  -2.83%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.773.59'
$ws.Range("E3").Value = '  -1.24%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.72%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9999'
$ws.Range("E5").Value = '  -0.82%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.04'
$ws.Range("E6").Value = '  -1.53%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4224'
$ws.Range("E7").Value = '  +0.85%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3587'
$ws.Range("E8").Value = '  +0.53%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07136'
$ws.Range("E9").Value = '  +0.83%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8348'
$ws.Range("E10").Value = '  -1.13%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.36'
$ws.Range("E11").Value = '  +1.28%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.771.56'
$ws.Range("E12").Value = '  -2.10%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.435'
$ws.Range("E13").Value = '  +1.40%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.236'
$ws.Range("E14").Value = '  -0.69%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06862'
$ws.Range("E15").Value = '  +0.15%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.96%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.94'
$ws.Range("E17").Value = '  -1.30%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008631'
$ws.Range("E18").Value = '  -1.14%  '

# Row 19
$ws.Range("E19").Value = '  -0.74%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.86'
$ws.Range("E20").Value = '  -1.29%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.350.20'
$ws.Range("E21").Value = '  -3.11%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.068'
$ws.Range("E22").Value = '  +0.45%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.89'
$ws.Range("E23").Value = '  +1.69%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.984.63'
$ws.Range("E24").Value = '  -3.05%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.64'
$ws.Range("E25").Value = '  -0.95%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.796'
$ws.Range("E26").Value = '  -8.42%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.94'
$ws.Range("E27").Value = '  -1.23%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.058'
$ws.Range("E28").Value = '  +0.88%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.37'
$ws.Range("E29").Value = '  +1.42%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.822'
$ws.Range("E30").Value = '  +9.34%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08813'
$ws.Range("E31").Value = '  -0.91%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7242'
$ws.Range("E32").Value = '  -0.24%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.118'
$ws.Range("E33").Value = '  +3.55%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.309'
$ws.Range("E34").Value = '  -1.15%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9991'
$ws.Range("E35").Value = '  -0.94%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.731'
$ws.Range("E36").Value = '  -5.60%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.083'
$ws.Range("E37").Value = '  +1.00%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05119'
$ws.Range("E38").Value = '  +0.11%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01877'
$ws.Range("E39").Value = '  -0.93%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4896'
$ws.Range("E40").Value = '  -1.21%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1602'
$ws.Range("E41").Value = '  -1.03%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.594'
$ws.Range("E42").Value = '  -2.73%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.294'
$ws.Range("E43").Value = '  +1.96%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.944'
$ws.Range("E44").Value = '  -1.13%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.21'
$ws.Range("E45").Value = '  +0.11%  '

# Row 46
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.22'
$ws.Range("E46").Value = '  -0.56%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9990'
$ws.Range("E47").Value = '  -0.98%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.627'
$ws.Range("E48").Value = '  +2.31%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06169'
$ws.Range("E49").Value = '  -2.28%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4438'
$ws.Range("E50").Value = '  -2.17%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.716'
$ws.Range("E51").Value = '  +2.69%  '
